# Add season record columns (Wins, Losses, Ties) to the PHI_2019 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold, bordered) from an existing header cell (AC1) to the new headers
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Re-set the values after paste (paste formats only shouldn't touch values, but just to be safe)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record data for each player row (2-58)
for ($row = 2; $row -le 58; $row++) {
    $ws.Cells.Item($row, 30).Value = 81  # AD = col 30 -> Wins
    $ws.Cells.Item($row, 31).Value = 81  # AE = col 31 -> Losses
    $ws.Cells.Item($row, 32).Value = 0   # AF = col 32 -> Ties
}
